{"js": "// Replace the division-problem text in each table cell with the new\n// problem, preserving all run/paragraph formatting (the cells are matched\n// and updated in place via Range.insertText, which keeps the existing\n// run properties).\nconst replacements = [\n  [\"811\u00f74=\", \"742\u00f78=\"],\n  [\"159\u00f75=\", \"241\u00f72=\"],\n  [\"426\u00f72=\", \"590\u00f72=\"],\n  [\"527\u00f78=\", \"343\u00f73=\"],\n  [\"846\u00f76=\", \"606\u00f76=\"],\n  [\"914\u00f79=\", \"182\u00f76=\"],\n  [\"858\u00f72=\", \"628\u00f72=\"],\n  [\"585\u00f79=\", \"609\u00f75=\"],\n  [\"803\u00f75=\", \"888\u00f77=\"],\n  [\"532\u00f72=\", \"975\u00f79=\"],\n  [\"535\u00f74=\", \"454\u00f79=\"],\n  [\"163\u00f75=\", \"515\u00f78=\"],\n  [\"447\u00f73=\", \"417\u00f76=\"],\n  [\"812\u00f76=\", \"397\u00f76=\"],\n  [\"653\u00f79=\", \"992\u00f78=\"],\n  [\"145\u00f72=\", \"102\u00f74=\"],\n  [\"356\u00f79=\", \"143\u00f74=\"],\n  [\"710\u00f73=\", \"809\u00f74=\"],\n  [\"314\u00f75=\", \"798\u00f78=\"],\n  [\"378\u00f79=\", \"312\u00f75=\"],\n  [\"349\u00f79=\", \"864\u00f74=\"],\n  [\"891\u00f76=\", \"135\u00f79=\"],\n  [\"645\u00f76=\", \"508\u00f73=\"],\n  [\"934\u00f79=\", \"649\u00f74=\"],\n  [\"394\u00f74=\", \"312\u00f72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each division-problem text in the table with its new value.\n# wdFindContinue = 1, wdReplaceAll = 2\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"811\u00f74=\", \"742\u00f78=\"),\n    @(\"159\u00f75=\", \"241\u00f72=\"),\n    @(\"426\u00f72=\", \"590\u00f72=\"),\n    @(\"527\u00f78=\", \"343\u00f73=\"),\n    @(\"846\u00f76=\", \"606\u00f76=\"),\n    @(\"914\u00f79=\", \"182\u00f76=\"),\n    @(\"858\u00f72=\", \"628\u00f72=\"),\n    @(\"585\u00f79=\", \"609\u00f75=\"),\n    @(\"803\u00f75=\", \"888\u00f77=\"),\n    @(\"532\u00f72=\", \"975\u00f79=\"),\n    @(\"535\u00f74=\", \"454\u00f79=\"),\n    @(\"163\u00f75=\", \"515\u00f78=\"),\n    @(\"447\u00f73=\", \"417\u00f76=\"),\n    @(\"812\u00f76=\", \"397\u00f76=\"),\n    @(\"653\u00f79=\", \"992\u00f78=\"),\n    @(\"145\u00f72=\", \"102\u00f74=\"),\n    @(\"356\u00f79=\", \"143\u00f74=\"),\n    @(\"710\u00f73=\", \"809\u00f74=\"),\n    @(\"314\u00f75=\", \"798\u00f78=\"),\n    @(\"378\u00f79=\", \"312\u00f75=\"),\n    @(\"349\u00f79=\", \"864\u00f74=\"),\n    @(\"891\u00f76=\", \"135\u00f79=\"),\n    @(\"645\u00f76=\", \"508\u00f73=\"),\n    @(\"934\u00f79=\", \"649\u00f74=\"),\n    @(\"394\u00f74=\", \"312\u00f72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
